$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9, shifting rows 9:23 down to 10:24.
$ws.Rows.Item(9).Insert()

# The new row 9 keeps the same "template" values as the row that used to be
# at position 9 (mercado/region/categoria/calidad/unidad/origen/kg-o-unidades/
# clasificacion), but gets a fresh date + volumen/precio data point.
$ws.Cells.Item(9, 1).Value = 1
$ws.Cells.Item(9, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(9, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(9, 4).Value = Get-Date -Year 2023 -Month 4 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(9, 5).Value = 15
$ws.Cells.Item(9, 6).Value = 100112001
$ws.Cells.Item(9, 7).Value = "Berenjena"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 190
$ws.Cells.Item(9, 11).Value = 4000
$ws.Cells.Item(9, 12).Value = 5000
$ws.Cells.Item(9, 13).Value = 4526
$ws.Cells.Item(9, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(9, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(9, 16).Value = 75
$ws.Cells.Item(9, 17).Value = 60
$ws.Cells.Item(9, 18).Value = "Hortaliza"
